$wb = $excel.ActiveWorkbook

function Set-DateCell($ws, $row, $value) {
    # Sets column-A cell to a numeric Excel date-serial value formatted as a date/time string
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $value
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

function Set-DataRow($ws, $row, $a, $aIsText, $b, $c, $d, $e, $f, $g, $h, $i) {
    # $aIsText: $true -> column A stays a plain text timestamp string (unconverted / "live" row)
    #           $false -> column A is a numeric date-serial value with date NumberFormat
    if ($aIsText) {
        $ws.Cells.Item($row, 1).Value = $a
    } else {
        Set-DateCell $ws $row $a
    }
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

# ==== ROW50-FE-LIFTER ====
$ws = $wb.Worksheets.Item("ROW50-FE-LIFTER")

# Convert existing inline-string date cells (col A) to numeric date-serial cells
Set-DateCell $ws 65 45721.73015877315
Set-DateCell $ws 66 45721.73018203703
Set-DateCell $ws 67 45721.7302053125
Set-DateCell $ws 68 45722.23047579861
Set-DateCell $ws 69 45722.23049802084
Set-DateCell $ws 70 45722.23052140047
Set-DateCell $ws 71 45723.19127907408
Set-DateCell $ws 72 45723.19130241898
Set-DateCell $ws 73 45723.19132570602

# Append new rows
Set-DataRow $ws 74 "2025-03-07 16:35:30" $True "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c," "0x01,0x90," "0x14" 400 568631262647113769549824.0 400 20

# ==== ROW50-MID-LIFTER ====
$ws = $wb.Worksheets.Item("ROW50-MID-LIFTER")

# Convert existing inline-string date cells (col A) to numeric date-serial cells
Set-DateCell $ws 65 45721.72839921296
Set-DateCell $ws 66 45721.7284225
Set-DateCell $ws 67 45721.72844564815
Set-DateCell $ws 68 45722.22854285879
Set-DateCell $ws 69 45722.22856443287
Set-DateCell $ws 70 45722.22858758102
Set-DateCell $ws 71 45722.72868413194
Set-DateCell $ws 72 45722.72870657407
Set-DateCell $ws 73 45722.72872986111

# Append new rows
Set-DataRow $ws 74 45723.22882704861 $False "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 568631262647113769549824.0 400 25
Set-DataRow $ws 75 45723.22884876157 $False "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 568631262647113769549824.0 400 25
Set-DataRow $ws 76 45723.22887202547 $False "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 568631262647113769549824.0 400 25
Set-DataRow $ws 77 "2025-03-07 17:29:34" $True "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x90," "0x19" 400 568631262647113769549824.0 400 25

# ==== ROW11-FE-LIFTER ====
$ws = $wb.Worksheets.Item("ROW11-FE-LIFTER")

# Convert existing inline-string date cells (col A) to numeric date-serial cells
Set-DateCell $ws 68 45721.72979140046
Set-DateCell $ws 69 45721.72981465278
Set-DateCell $ws 70 45721.72983799769
Set-DateCell $ws 71 45722.23010865741
Set-DateCell $ws 72 45722.23013063658
Set-DateCell $ws 73 45722.23015388889
Set-DateCell $ws 74 45723.19126707176
Set-DateCell $ws 75 45723.19129023148

# Append new rows
Set-DataRow $ws 76 45723.19131362269 $False "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x90," "0x13" 400 568631262647113769549824.0 400 19
Set-DataRow $ws 77 "2025-03-07 16:35:29" $True "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x90," "0x13" 400 568631262647113769549824.0 400 19

# ==== ROW11-MID-LIFTER ====
$ws = $wb.Worksheets.Item("ROW11-MID-LIFTER")

# Convert existing inline-string date cells (col A) to numeric date-serial cells
Set-DateCell $ws 65 45721.73010056713
Set-DateCell $ws 66 45721.73012371528
Set-DateCell $ws 67 45721.73014709491
Set-DateCell $ws 68 45722.23024466435
Set-DateCell $ws 69 45722.23026591435
Set-DateCell $ws 70 45722.23028918981
Set-DateCell $ws 71 45723.19114016204
Set-DateCell $ws 72 45723.19116321759
Set-DateCell $ws 73 45723.19118659722

# Append new rows
Set-DataRow $ws 74 "2025-03-07 16:35:18" $True "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x90," "0x9" 400 568631262647113769549824.0 400 9
